$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.200453400611877
$ws.Range("B1").Value = 1.666729927062988
$ws.Range("C1").Value = 3.649246215820312
$ws.Range("D1").Value = 3.440254926681519
$ws.Range("E1").Value = 0.9707048535346985
